# Update the cryptocurrency price/volume symbol list to the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new value (as text, to preserve exact formatting)
$updates = @{
    "D2"  = "244.83"
    "D3"  = "23.09"
    "D5"  = "0.05982"
    "D7"  = "0.8087"
    "D8"  = "0.9276"
    "D9"  = "0.1424"
    "D10" = "0.07442"
    "D12" = "0.03040"
    "D13" = "0.09355"
    "D14" = "3.939"
    "D15" = "0.001587"
    "D16" = "0.04820"
    "E17" = "16OneONEWorstin24h"
    "D18" = "0.005678"
    "D20" = "0.0009869"
    "D22" = "3.658"
    "D23" = "6.451"
    "D40" = "0.03941"
    "D41" = "0.006182"
    "E41" = "40KickTokenKICK"
    "D44" = "0.007535"
    "E44" = "43LocalTradersLCTBestin24h"
    "D45" = "0.00005209"
    "E49" = "48BOLOBOLO"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}
